$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =====================================================================
# Row 48: the earlier "IN PROGRESS" Sell of 209 XRP is now finalized.
# Status flips to DONE and the Finalized date / Fee / Profit / Duration
# columns (I:L) get filled in for the first time.
# =====================================================================

# I48 - Finalized date (numeric date, must land on the same style=6 the
# empty placeholder cell already carries).
$finDate = Get-Date -Year 2017 -Month 5 -Day 8 -Hour 9 -Minute 33 -Second 59
$ws.Range("I48").Value = $finDate

# J48 - Fee
$ws.Range("J48").Value = "0.06196192 USDT (0.15%)"

# L48 - Transaction duration
$ws.Range("L48").Value = "3 day"

# K48 - Profit(%), rendered as "     " + green "~6%" (matches the other
# Profit cells in the sheet).
$ws.Range("K48").Value = "     ~6%"
$ws.Range("K48").Characters(6, 3).Font.Color = 5287936

# H48 - Status DONE (set last so the earlier placeholder text doesn't
# linger if anything above fails)
$ws.Range("H48").Value = "DONE"

# =====================================================================
# Row 49 (new row): the next trade - a Buy of 225 XRP - freshly logged
# and still "IN PROGRESS".
# =====================================================================

# A49 - Data (date), same instant the previous row got finalized.
$ws.Range("A49").Value = $finDate
$ws.Range("A47").Copy()
$ws.Range("A49").PasteSpecial(-4122)
$ws.Range("A49").Value = $finDate

# B49 - Action(Buy/Sell): "Buy" in green
$ws.Range("B49").Value = "            Buy"
$ws.Range("B49").Characters(13, 3).Font.Color = 5287936

# C49 - Currency
$ws.Range("C49").Value = "        XRP"

# D49 - Data (wrapped text that looks numeric, so it must be entered
# with a leading quote to keep it text, then have the quote-prefix
# style swapped out for the plain wrap-text style already used by
# D45:D48).
$ws.Range("D49").Value = "'           0.19130000`n"
$ws.Range("D47").Copy()
$ws.Range("D49").PasteSpecial(-4122)

# E49 - Transaction value
$ws.Range("E49").Value = "         0.185USDT"

# F49 - Transaction amount
$ws.Range("F49").Value = "         225 XRP"

# G49 - Transaction code
$ws.Range("G49").Value = " XRP/USDT0000005"

# H49 - Status
$ws.Range("H49").Value = "IN PROGRESS"

# I49 - Finalized date placeholder (still empty, but keep the date
# style so the column stays consistent).
$ws.Range("A45").Copy()
$ws.Range("I49").PasteSpecial(-4122)
$ws.Range("I49").ClearContents()

# K49 - Profit(%) placeholder (blank spaces, no color yet).
$ws.Range("K49").Value = "     "

$ws.Rows.Item(49).RowHeight = 14.25

$ws.Range("K49").Select()
